$d = $word.ActiveDocument

# The document ends with a paragraph whose run is followed by the
# Word-managed "_GoBack" bookmark (marking the last edit position).
# We need to add a brand-new paragraph ("我是智豪") right after the
# existing last paragraph, and have the _GoBack bookmark end up at the
# end of that new paragraph (i.e. move along with the new "last edit").

# 1) Remove the existing _GoBack bookmark from its current location.
#    (It is a hidden bookmark, so it won't show up in $d.Bookmarks'
#    enumeration/Count, but it can still be reached by name.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

# 2) Find the insertion point: right at the end of the text of the
#    current last paragraph (i.e. just before its paragraph mark).
$lastPara = $d.Paragraphs.Last
$insertPos = $lastPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

# 3) Insert a whole new paragraph, with the new run of text and the
#    _GoBack bookmark re-created at its end, via a raw WordprocessingML
#    fragment (keeps the bookmark placement exact and avoids any
#    ambiguity around paragraph-boundary Range positions).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = '<w:p ' + $wNs + '>' +
    '<w:pPr>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
        '<w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr>' +
        '<w:t>我是智豪</w:t>' +
    '</w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

[void]$insertRange.InsertXML($newParaXml)
